$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.833.86'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '2.924.62'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.27%  '
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0883'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.136'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").Value = '3.382.25'
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("D16").Value = '2.913.96'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.979'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("D18").Value = '51.835.38'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("E19").Value = '  -1.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("D22").Value = '0.0₃0982'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("E26").Value = '  +10.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("E29").Value = '  +13.44%  '
$ws.Range("E30").Value = '  +11.98%  '
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '52.19'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("E35").Value = '  -4.36%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -14.98%  '
$ws.Range("E38").Value = '  -2.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.42'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.93%  '
$ws.Range("E41").Value = '  +4.95%  '
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '120.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.69%  '
$ws.Range("E47").Value = '  -4.40%  '
$ws.Range("D48").Value = '2.139.75'
$ws.Range("E48").Value = '  -2.94%  '
$ws.Range("E49").Value = '  -6.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0339'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.05%  '
$ws.Range("E51").Value = '  -4.25%  '
